$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 1052
$ws.Cells.Item(3, 6).Value = 658
$ws.Cells.Item(4, 6).Value = 1460
$ws.Cells.Item(6, 6).Value = 3203
$ws.Cells.Item(8, 6).Value = 596
$ws.Cells.Item(9, 6).Value = 2154
$ws.Cells.Item(10, 6).Value = 459
$ws.Cells.Item(11, 6).Value = 389
$ws.Cells.Item(13, 6).Value = 117
$ws.Cells.Item(14, 6).Value = 271
$ws.Cells.Item(16, 6).Value = 1047
$ws.Cells.Item(17, 6).Value = 419
$ws.Cells.Item(18, 6).Value = 71
$ws.Cells.Item(19, 6).Value = 176
$ws.Cells.Item(20, 6).Value = 4250
$ws.Cells.Item(21, 6).Value = 1248
$ws.Cells.Item(22, 6).Value = 3287
$ws.Cells.Item(24, 6).Value = 126
$ws.Cells.Item(25, 6).Value = 3135
$ws.Cells.Item(26, 6).Value = 4761
$ws.Cells.Item(28, 6).Value = 963
$ws.Cells.Item(29, 6).Value = 525
$ws.Cells.Item(30, 6).Value = 3098
$ws.Cells.Item(31, 6).Value = 314
$ws.Cells.Item(33, 6).Value = 123
$ws.Cells.Item(36, 6).Value = 1120
$ws.Cells.Item(37, 6).Value = 1374
$ws.Cells.Item(39, 6).Value = 1273
$ws.Cells.Item(40, 6).Value = 815
$ws.Cells.Item(42, 6).Value = 760
$ws.Cells.Item(43, 6).Value = 483
$ws.Cells.Item(44, 6).Value = 46
$ws.Cells.Item(45, 6).Value = 253
$ws.Cells.Item(46, 6).Value = 53
$ws.Cells.Item(47, 6).Value = 112
$ws.Cells.Item(49, 6).Value = 3690

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(6, 6).Value = 981

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 1959

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 1959
$ws.Cells.Item(4, 6).Value = 658
$ws.Cells.Item(5, 6).Value = 1460
$ws.Cells.Item(6, 6).Value = 3203
$ws.Cells.Item(8, 6).Value = 2154
$ws.Cells.Item(9, 6).Value = 459
$ws.Cells.Item(10, 6).Value = 389
$ws.Cells.Item(13, 6).Value = 981
$ws.Cells.Item(15, 6).Value = 117
$ws.Cells.Item(16, 6).Value = 271
$ws.Cells.Item(17, 6).Value = 1047
$ws.Cells.Item(19, 6).Value = 419
$ws.Cells.Item(20, 6).Value = 176
$ws.Cells.Item(21, 6).Value = 4250
$ws.Cells.Item(23, 6).Value = 1248
$ws.Cells.Item(25, 6).Value = 3287
$ws.Cells.Item(26, 6).Value = 3136
$ws.Cells.Item(27, 6).Value = 4761
$ws.Cells.Item(28, 6).Value = 963
$ws.Cells.Item(29, 6).Value = 3099
$ws.Cells.Item(30, 6).Value = 314
$ws.Cells.Item(32, 6).Value = 123
$ws.Cells.Item(35, 6).Value = 1121
$ws.Cells.Item(36, 6).Value = 1374
$ws.Cells.Item(38, 6).Value = 1273
$ws.Cells.Item(39, 6).Value = 815
$ws.Cells.Item(41, 6).Value = 483
$ws.Cells.Item(43, 6).Value = 46
$ws.Cells.Item(45, 6).Value = 253
$ws.Cells.Item(46, 6).Value = 53
$ws.Cells.Item(47, 6).Value = 112
$ws.Cells.Item(49, 6).Value = 3690
